$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45008
$ws.Range("D4").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101004
$ws.Range("J4").Value = "Frambuesa"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("Q4").Value = '$/bandeja 2 kilos'
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value = 3500
$ws.Range("T4").Value = 2

# Row 5
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 45008
$ws.Range("D5").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101004
$ws.Range("J5").Value = "Frambuesa"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("Q5").Value = '$/bandeja 2 kilos'
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 2
